$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.129414558410645
$ws.Range("B1").Value = 3.591049909591675
$ws.Range("C1").Value = 4.165102958679199
$ws.Range("D1").Value = 2.696932315826416
$ws.Range("E1").Value = 1.052422285079956
